# Scheduled market-data refresh: recompute crafting-profit columns
# (currentAveragePrice / *NQ / *HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ -> cols H:N)
# for the leves whose linked-item market snapshot changed since the last run.
# Values below are the new Universalis-sourced snapshot for each (sheet,row).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 21: "Book and a Hard Place" / Engraved Hard Leather Grimoire
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()

# Row 23: "There's Something about Bury" / Hard Leather Grimoire
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()

# Row 29: "Dripping with Venom" / Weak Blinding Potion
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("N29").ClearContents()

# Row 38: "Just Give Him a Serum" / Hi-Potion of Strength
$ws.Range("H38").Value = 676.7273
$ws.Range("I38").Value = 171.55556
$ws.Range("J38").Value = 2950
$ws.Range("K38").Value = 514.66668
$ws.Range("L38").Value = 8850
$ws.Range("M38").Value = -142.66668
$ws.Range("N38").Value = -9594

# Row 58: "A Matter of Vital Importance" / Mega-Potion of Vitality
$ws.Range("H58").Value = 978.3570999999999
$ws.Range("I58").Value = 86.666664
$ws.Range("J58").Value = 2583.4
$ws.Range("K58").Value = 259.999992
$ws.Range("L58").Value = 7750.200000000001
$ws.Range("M58").Value = -109.999992
$ws.Range("N58").Value = -8050.200000000001

# Row 64: "Forged from the Void" / Void Glue
$ws.Range("H64").Value = 641405.0600000001
$ws.Range("I64").Value = 1023588.1
$ws.Range("J64").Value = 4433.3335
$ws.Range("K64").Value = 1023588.1
$ws.Range("L64").Value = 4433.3335
$ws.Range("M64").Value = -1023340.1
$ws.Range("N64").Value = -4929.3335

# Row 67: "Dodging the Draft (L)" / Void Glue
$ws.Range("H67").Value = 641405.0600000001
$ws.Range("I67").Value = 1023588.1
$ws.Range("J67").Value = 4433.3335
$ws.Range("K67").Value = 1023588.1
$ws.Range("L67").Value = 4433.3335
$ws.Range("M67").Value = -1022730.1
$ws.Range("N67").Value = -6149.3335

# Row 69: "Steeling the Knife, Steeling the Mind" / Grade 1 Mind Dissolvent
$ws.Range("H69").Value = 4507.5
$ws.Range("J69").Value = 3409
$ws.Range("L69").Value = 10227
$ws.Range("N69").Value = -11975

# Row 72: "Surgical Substitution (L)" / Grade 1 Mind Dissolvent
$ws.Range("H72").Value = 4507.5
$ws.Range("J72").Value = 3409
$ws.Range("L72").Value = 30681
$ws.Range("N72").Value = -39417

$ws = $wb.Worksheets.Item("ARM")
# Row 10: "Bronzed and Burnt" / Bronze Sallet
$ws.Range("H10").Value = 850
$ws.Range("I10").Value = 850
$ws.Range("K10").Value = 850
$ws.Range("M10").Value = -680

# Row 61: "Dealing with the Tough Stuff" / Cobalt Ingot
$ws.Range("H61").Value = 9296.795
$ws.Range("I61").Value = 5902.231
$ws.Range("J61").Value = 16085.923
$ws.Range("K61").Value = 5902.231
$ws.Range("L61").Value = 16085.923
$ws.Range("M61").Value = -5690.231
$ws.Range("N61").Value = -16509.923

# Row 74: "As the Bolt Flies" / Titanium Nugget
$ws.Range("H74").Value = 118803.5
$ws.Range("I74").Value = 142148.44
$ws.Range("J74").Value = 2078.8
$ws.Range("K74").Value = 142148.44
$ws.Range("L74").Value = 2078.8
$ws.Range("M74").Value = -141274.44
$ws.Range("N74").Value = -3826.8

# Row 77: "Heavy Metal Banned (L)" / Titanium Nugget
$ws.Range("H77").Value = 118803.5
$ws.Range("I77").Value = 142148.44
$ws.Range("J77").Value = 2078.8
$ws.Range("K77").Value = 710742.2
$ws.Range("L77").Value = 10394
$ws.Range("M77").Value = -706374.2
$ws.Range("N77").Value = -19130

# Row 136: "Metal with Mettle" / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 9296.795
$ws.Range("I136").Value = 5902.231
$ws.Range("J136").Value = 16085.923
$ws.Range("K136").Value = 17706.693
$ws.Range("L136").Value = 48257.769
$ws.Range("M136").Value = -15156.693
$ws.Range("N136").Value = -53357.769

$ws = $wb.Worksheets.Item("CRP")
# Row 12: "A Sword in Hand" / Ash Macuahuitl
$ws.Range("H12").Value = 2166.6667
$ws.Range("I12").Value = 1000
$ws.Range("J12").Value = 2750
$ws.Range("K12").Value = 1000
$ws.Range("L12").Value = 2750
$ws.Range("M12").Value = -830
$ws.Range("N12").Value = -3090

# Row 110: "A Stronger Offense" / Applewood Spear
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 3: "Needful Rings" / Copper Wristlets
$ws.Range("H3").Value = 105101944
$ws.Range("J3").Value = 2780.4
$ws.Range("L3").Value = 2780.4
$ws.Range("N3").Value = -3012.4

# Row 12: "Horn of Plenty" / Bone Armillae
$ws.Range("H12").Value = 401.66666
$ws.Range("I12").Value = 401.66666
$ws.Range("K12").Value = 401.66666
$ws.Range("M12").Value = -261.66666

# Row 122: "Awarding Academic Excellence" / Ametrine
$ws.Range("H122").Value = 7712.375
$ws.Range("I122").Value = 17299.666
$ws.Range("J122").Value = 1960
$ws.Range("K122").Value = 51898.99800000001
$ws.Range("L122").Value = 5880
$ws.Range("M122").Value = -49448.99800000001
$ws.Range("N122").Value = -10780

# Row 126: "Gold Rush Order" / Phrygian Gold Ingot
$ws.Range("H126").Value = 2035.3429
$ws.Range("I126").Value = 1629.0454
$ws.Range("J126").Value = 2722.923
$ws.Range("K126").Value = 4887.1362
$ws.Range("L126").Value = 8168.768999999999
$ws.Range("M126").Value = -2417.1362
$ws.Range("N126").Value = -13108.769

# Row 132: "On Board for Lar" / Lar Ingot
$ws.Range("H132").Value = 12157.75
$ws.Range("I132").Value = 3922.625
$ws.Range("J132").Value = 28628
$ws.Range("K132").Value = 11767.875
$ws.Range("L132").Value = 85884
$ws.Range("M132").Value = -9237.875
$ws.Range("N132").Value = -90944

$ws = $wb.Worksheets.Item("LTW")
# Row 10: "In Their Shoes" / Leather Leggings
$ws.Range("H10").Value = 2085
$ws.Range("I10").Value = 868.6667
$ws.Range("J10").Value = 3301.3333
$ws.Range("K10").Value = 868.6667
$ws.Range("L10").Value = 3301.3333
$ws.Range("M10").Value = -728.6667
$ws.Range("N10").Value = -3581.3333

# Row 132: "Tenets of Tanning" / Silver Lobo Leather
$ws.Range("H132").Value = 3447.2778
$ws.Range("I132").Value = 3505.0476
$ws.Range("K132").Value = 10515.1428
$ws.Range("M132").Value = -7985.1428

# Row 136: "Respect for Br'aax" / Br'aax Leather
$ws.Range("H136").Value = 4249.6665
$ws.Range("I136").Value = 2642.484
$ws.Range("K136").Value = 7927.451999999999
$ws.Range("M136").Value = -5377.451999999999

$ws = $wb.Worksheets.Item("WVR")
# Row 14: "Hat in Hand" / Straw Hat
$ws.Range("H14").Value = 626350
$ws.Range("I14").Value = 2500000
$ws.Range("J14").Value = 1800
$ws.Range("K14").Value = 2500000
$ws.Range("L14").Value = 1800
$ws.Range("M14").Value = -2499832
$ws.Range("N14").Value = -2136

# Row 45: "Private Concerns" / Linen Trousers
$ws.Range("H45").Value = 4950
$ws.Range("J45").Value = 4950
$ws.Range("L45").Value = 4950
$ws.Range("M45").Value = 4950
$ws.Range("N45").Value = -5932

# Row 62: "Pride Up in Smoke" / Rainbow Cloth
$ws.Range("H62").Value = 3961.2
$ws.Range("J62").Value = 3944.5715
$ws.Range("L62").Value = 3944.5715
$ws.Range("N62").Value = -5192.5715

# Row 65: "Desperate for Diversionaries (L)" / Rainbow Cloth
$ws.Range("H65").Value = 3961.2
$ws.Range("J65").Value = 3944.5715
$ws.Range("L65").Value = 19722.8575
$ws.Range("N65").Value = -25962.8575

# Row 68: "What Not to Wear" / Holy Rainbow Shirt of Striking
$ws.Range("H68").Value = 37180.668
$ws.Range("J68").Value = 37180.668
$ws.Range("L68").Value = 37180.668
$ws.Range("N68").Value = -38802.668

# Row 71: "Appeal of Foreign Apparel (L)" / Holy Rainbow Shirt of Striking
$ws.Range("H71").Value = 37180.668
$ws.Range("J71").Value = 37180.668
$ws.Range("L71").Value = 111542.004
$ws.Range("N71").Value = -119654.004

# Row 81: "Where the Dragonflies, the Net Catches" / Crawler Silk
$ws.Range("H81").Value = 15386546
$ws.Range("I81").Value = 1357.4286
$ws.Range("J81").Value = 33335934
$ws.Range("K81").Value = 2714.8572
$ws.Range("L81").Value = 66671868
$ws.Range("M81").Value = -1653.8572
$ws.Range("N81").Value = -66673990

# Row 84: "To Kill a Dragon on Nameday (L)" / Crawler Silk
$ws.Range("H84").Value = 15386546
$ws.Range("I84").Value = 1357.4286
$ws.Range("J84").Value = 33335934
$ws.Range("K84").Value = 13574.286
$ws.Range("L84").Value = 333359340
$ws.Range("M84").Value = -8270.286
$ws.Range("N84").Value = -333369948
